# [PHOENIX-5854] - update grievancesTestData: replace the "Dog menace" /
# Public Health and Sanitation sample grievance with a "Street Lighting"
# scenario on the grievanceDetails sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grievanceDetails")

# Order matters: it reproduces the shared-string insertion order seen in
# the saved workbook (grievanceType, then grievanceDetails, then
# grievanceCategory).
$ws.Range("C2").Value = "Non Burning of Street Lights"
$ws.Range("D2").Value = "No street light past one week"
$ws.Range("B2").Value = "Street Lighting"

# Column E (grievanceLocation) widened to fit the new data.
$ws.Columns.Item(5).ColumnWidth = 31.66

# Leave the cursor parked on B5, matching the resaved workbook's view state.
$ws.Range("B5").Select() | Out-Null
